$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E25").Value = 126.78
$ws.Range("E26").Value = 138.63
